$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.17
$ws.Range("D3").Value = 0.08
$ws.Range("D5").Value = 0.07000000000000001
$ws.Range("D6").Value = 0.08
$ws.Range("D7").Value = 0.1
$ws.Range("D8").Value = 0.12
$ws.Range("D9").Value = 0.2
$ws.Range("D10").Value = 0.08
$ws.Range("D12").Value = 0.04
$ws.Range("D16").Value = 0.08
$ws.Range("D17").Value = 0.06
$ws.Range("D18").Value = 0.11
$ws.Range("D19").Value = 0.15
$ws.Range("D21").Value = 0.2
$ws.Range("D22").Value = 0.09
$ws.Range("D28").Value = 0.07000000000000001
$ws.Range("D29").Value = 0.12
$ws.Range("D30").Value = 0.09
$ws.Range("D31").Value = 0.11
$ws.Range("D33").Value = 0.11
$ws.Range("D36").Value = 0.08
